$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the new
# header cells so the new columns pick up the same bold/border/centered
# header style (style index 1) instead of minting a brand new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) gets the same team record: 83 wins, 79 losses, 0 ties.
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 83
    $ws.Cells.Item($row, 31).Value = 79
    $ws.Cells.Item($row, 32).Value = 0
}
